$wb = $excel.ActiveWorkbook

# --- Data updates on the "Training Dashboard" sheet ---
$ws = $wb.Worksheets.Item("Training Dashboard")

# H3: 92 -> 84
$ws.Range("H3").Value = 84

# I3: "08-Sep-2025" -> "16-Sep-2025" (kept as literal text, not an Excel date serial)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
# Restore the original (non-text) number format / style by copying the format
# from an untouched neighboring cell in the same row, then re-apply the value.
$ws.Range("J3").Copy()
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Header / title font re-styling (applies to both sheets) ---
# Title cell (A1) and the column-header row both move to a bold WHITE font
# (previously the title was bold/size-14 and the header was bold/black).
foreach ($sheetName in @("Training Dashboard", "Exam Dashboard")) {
    $sheet = $wb.Worksheets.Item($sheetName)

    $titleCell = $sheet.Range("A1")
    $titleCell.Font.Size = 11
    $titleCell.Font.Color = 16777215

    $headerRow = $sheet.Range("A2").EntireRow
    $usedHeader = $sheet.Range($sheet.Cells.Item(2, 1), $sheet.Cells.Item(2, 11))
    $usedHeader.Font.Size = 11
    $usedHeader.Font.Color = 16777215
}
